$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D"  = -0.0815
    "G"  = -0.3610619469026549
    "H"  = -0.3610619469026549
    "I"  = -0.3008849557522124
    "J"  = -0.3008849557522124
    "K"  = -1.7
    "L"  = -0.3008849557522124
    "U"  = 0.083
    "V"  = 0.1194244604316547
    "W"  = 0.6439393939393939
    "X"  = 0.2321979651145056
    "Y"  = 0.4117414288248883
    "Z"  = -2.10271678451805
    "AA" = 0.6326758466691477
    "AB" = 0.08368386324973645
    "AC" = 0.5489919834194112
    "AD" = 2.26
    "AF" = 2.26
    "AG" = 2.177
    "AH" = 0.7648054145516074
    "AI" = -1.018018018018018
    "AJ" = 0.758008356545961
    "AK" = -0.9452887537993916
    "AL" = 0.039
    "AM" = -0.019
    "AN" = -1.614285714285714
    "AO" = -43.58974358974359
    "AP" = -1.555
    "AQ" = 89.4736842105263
}

foreach ($row in 2..3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
